# Auto-generated Excel COM-interop script
# Applies a data refresh to the Louisoix_Profits workbook per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2835.8
$ws.Range("I40").Value = 1475.5
$ws.Range("J40").Value = 3045.077
$ws.Range("K40").Value = 1475.5
$ws.Range("L40").Value = 3045.077
$ws.Range("M40").Value = -1300.5
$ws.Range("N40").Value = -3395.077

$ws.Range("H80").Value = 1465.75
$ws.Range("I80").Value = 2281.6667
$ws.Range("K80").Value = 6845.000100000001
$ws.Range("M80").Value = -5847.000100000001

$ws.Range("H83").Value = 1465.75
$ws.Range("I83").Value = 2281.6667
$ws.Range("K83").Value = 20535.0003
$ws.Range("M83").Value = -15543.0003

$ws.Range("H137").Value = 2337.652
$ws.Range("J137").Value = 2898.5386
$ws.Range("L137").Value = 8695.6158
$ws.Range("N137").Value = -13795.6158

$ws.Range("H141").Value = 1312.4
$ws.Range("I141").Value = 1347.6666
$ws.Range("J141").Value = 995
$ws.Range("K141").Value = 4042.9998
$ws.Range("L141").Value = 2985
$ws.Range("M141").Value = 1137.0002
$ws.Range("N141").Value = -13345

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4990.6
$ws.Range("I61").Value = 2385.2222
$ws.Range("K61").Value = 2385.2222
$ws.Range("M61").Value = -2173.2222

$ws.Range("H74").Value = 742.7143
$ws.Range("I74").Value = 549.75
$ws.Range("J74").Value = 1000
$ws.Range("K74").Value = 549.75
$ws.Range("L74").Value = 1000
$ws.Range("M74").Value = 324.25
$ws.Range("N74").Value = -2748

$ws.Range("H77").Value = 742.7143
$ws.Range("I77").Value = 549.75
$ws.Range("J77").Value = 1000
$ws.Range("K77").Value = 2748.75
$ws.Range("L77").Value = 5000
$ws.Range("M77").Value = 1619.25
$ws.Range("N77").Value = -13736

$ws.Range("H97").Value = 8244.5625
$ws.Range("I97").Value = 10675.9
$ws.Range("K97").Value = 10675.9
$ws.Range("M97").Value = -10179.9

$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws.Range("H136").Value = 4990.6
$ws.Range("I136").Value = 2385.2222
$ws.Range("K136").Value = 7155.6666
$ws.Range("M136").Value = -4605.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 167698.67
$ws.Range("J22").Value = 1201
$ws.Range("L22").Value = 1201
$ws.Range("N22").Value = -1547

$ws.Range("H103").Value = 40910.6
$ws.Range("J103").Value = 40910.6
$ws.Range("L103").Value = 40910.6
$ws.Range("N103").Value = -43254.6

$ws.Range("H134").Value = 2583.8647
$ws.Range("I134").Value = 2186.8823
$ws.Range("J134").Value = 7083
$ws.Range("K134").Value = 6560.646900000001
$ws.Range("L134").Value = 21249
$ws.Range("M134").Value = -4025.646900000001
$ws.Range("N134").Value = -26319

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 129.8
$ws.Range("I2").Value = 87.25
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 87.25
$ws.Range("L2").Value = 300
$ws.Range("M2").Value = 25.75
$ws.Range("N2").Value = -526

$ws.Range("H3").Value = 1301
$ws.Range("I3").Value = 602
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 602
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -489
$ws.Range("N3").Value = -2226

$ws.Range("H31").Value = 2633
$ws.Range("I31").Value = 1916.5
$ws.Range("J31").Value = 5499
$ws.Range("K31").Value = 1916.5
$ws.Range("L31").Value = 5499
$ws.Range("M31").Value = -1621.5
$ws.Range("N31").Value = -6089

$ws.Range("H34").Value = 2633
$ws.Range("I34").Value = 1916.5
$ws.Range("J34").Value = 5499
$ws.Range("K34").Value = 1916.5
$ws.Range("L34").Value = 5499
$ws.Range("M34").Value = -1714.5
$ws.Range("N34").Value = -5903

$ws.Range("H58").Value = 79685.92
$ws.Range("I58").Value = 102570.3
$ws.Range("K58").Value = 102570.3
$ws.Range("M58").Value = -102367.3

$ws.Range("H86").Value = 6841.2856
$ws.Range("I86").Value = 7473.75
$ws.Range("J86").Value = 5998
$ws.Range("K86").Value = 7473.75
$ws.Range("L86").Value = 5998
$ws.Range("M86").Value = -6350.75
$ws.Range("N86").Value = -8244

$ws.Range("H89").Value = 6841.2856
$ws.Range("I89").Value = 7473.75
$ws.Range("J89").Value = 5998
$ws.Range("K89").Value = 37368.75
$ws.Range("L89").Value = 29990
$ws.Range("M89").Value = -31752.75
$ws.Range("N89").Value = -41222

$ws.Range("H122").Value = 2541.8
$ws.Range("I122").Value = 2541.8
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7625.400000000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5175.400000000001
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 2551.1052
$ws.Range("I132").Value = 2071.4285
$ws.Range("K132").Value = 6214.2855
$ws.Range("M132").Value = -3684.2855

$ws.Range("H134").Value = 72543
$ws.Range("I134").Value = 107147.09
$ws.Range("K134").Value = 321441.27
$ws.Range("M134").Value = -318906.27

$ws.Range("H136").Value = 79685.92
$ws.Range("I136").Value = 102570.3
$ws.Range("K136").Value = 307710.9
$ws.Range("M136").Value = -305160.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 45066.6
$ws.Range("J11").Value = 166249.75
$ws.Range("L11").Value = 166249.75
$ws.Range("N11").Value = -166527.75

$ws.Range("H100").Value = 35300
$ws.Range("J100").Value = 35300
$ws.Range("L100").Value = 35300
$ws.Range("N100").Value = -37464

$ws.Range("H122").Value = 4368.125
$ws.Range("I122").Value = 2074.5
$ws.Range("J122").Value = 5132.6665
$ws.Range("K122").Value = 6223.5
$ws.Range("L122").Value = 15397.9995
$ws.Range("M122").Value = -3773.5
$ws.Range("N122").Value = -20297.9995

$ws.Range("H132").Value = 65296.25
$ws.Range("I132").Value = 69449.336
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 208348.008
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -205818.008
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5474.875
$ws.Range("I7").Value = 3399.8572
$ws.Range("K7").Value = 3399.8572
$ws.Range("M7").Value = -3287.8572

$ws.Range("H46").Value = 4367.4062
$ws.Range("I46").Value = 11818.818
$ws.Range("K46").Value = 11818.818
$ws.Range("M46").Value = -11630.818

$ws.Range("H93").Value = 2287.8333
$ws.Range("I93").Value = 2132.6365
$ws.Range("K93").Value = 2132.6365
$ws.Range("M93").Value = -884.6365000000001

$ws.Range("H126").Value = 5474.875
$ws.Range("I126").Value = 3399.8572
$ws.Range("K126").Value = 10199.5716
$ws.Range("M126").Value = -7729.571599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1835.64
$ws.Range("I81").Value = 1528.0476
$ws.Range("J81").Value = 3450.5
$ws.Range("K81").Value = 3056.0952
$ws.Range("L81").Value = 6901
$ws.Range("M81").Value = -1995.0952
$ws.Range("N81").Value = -9023

$ws.Range("H84").Value = 1835.64
$ws.Range("I84").Value = 1528.0476
$ws.Range("J84").Value = 3450.5
$ws.Range("K84").Value = 15280.476
$ws.Range("L84").Value = 34505
$ws.Range("M84").Value = -9976.476000000001
$ws.Range("N84").Value = -45113

$ws.Range("H132").Value = 87319.664
$ws.Range("I132").Value = 95076
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 285228
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -282698
$ws.Range("N132").Value = -11060

